# Fixed truncation error in effluent data table.
#
# The "effluent" sheet had a stray, truncated 3rd column ("conc") and
# placeholder "NA" text values left over from an earlier layout. This
# restores the intended 2-column table (hours / <ion name>) with real
# numeric data instead of "NA" placeholders, and drops the now-unused
# "conc" column so effluent output isn't silently cut off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effluent")

# Column B now holds the ion name ("CHLORIDE") instead of the generic
# "name" header.
$ws.Range("B1").Value = "CHLORIDE"

# Replace the "NA" placeholder row with real numeric data.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0

# Drop the truncated 3rd ("conc") column entirely.
$ws.Columns.Item(3).Delete()

# Restore the selection to match the saved view.
$ws.Range("L4").Select() | Out-Null
